$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "s"

# --- New player roster data (jersey, player, date, opponent, home, position) ---
$players = @(
  @(0,  "Kevin_Love",        "Boston_Celtics", 1, "F"),
  @(1,  "Derrick_Rose",      "Boston_Celtics", 1, "G"),
  @(3,  "isahiah_Thomas",    "Boston_Celtics", 1, "G"),
  @(4,  "Iman_Shumpert",     "Boston_Celtics", 1, "G"),
  @(5,  "Jr.Smith",          "Boston_Celtics", 1, "GF"),
  @(8,  "Channing Fyre",     "Boston_Celtics", 1, "F"),
  @(9,  "Dwyane _Wade",      "Boston_Celtics", 1, "G"),
  @(10, "Jhon_Holland",      "Boston_Celtics", 1, "GF"),
  @(13, "Tristan-Thompson",  "Boston_Celtics", 1, "CF"),
  @(16, "Cedi_Osman",        "Boston_Celtics", 1, "F"),
  @(23, "Lebron_James",      "Boston_Celtics", 1, "F"),
  @(26, "Kyle_Korver",       "Boston_Celtics", 1, "G"),
  @(32, "Jeff_Green",        "Boston_Celtics", 1, "F"),
  @(41, "Ante_Zizic",        "Boston_Celtics", 1, "FC"),
  @(81, "Jose_Calderon",     "Boston_Celtics", 1, "G"),
  @(99, "Jae_Crowder",       "Boston_Celtics", 1, "F")
)

$r = 2
foreach ($p in $players) {
  $ws.Cells.Item($r, 1).Value = $p[0]
  $ws.Cells.Item($r, 2).Value = $p[1]
  $ws.Cells.Item($r, 3).Value = "10/17/2017"
  $ws.Cells.Item($r, 4).Value = $p[2]
  $ws.Cells.Item($r, 5).Value = $p[3]
  $ws.Cells.Item($r, 6).Value = $p[4]
  $r = $r + 1
}

# --- Remove now-unused trailing rows (previously 18,19,20) ---
$ws.Rows("18:20").Delete()

# --- Apply shared style (font) to columns C:Y for header + data rows ---
$ws.Range("C1:Y17").Font.Name = "Calibri"
$ws.Range("C1:Y17").Font.Size = 11

# --- Apply date number format to column C (data rows only) ---
$ws.Range("C2:C17").NumberFormat = "mm-dd-yy"
